# The workbook now lets the script print links to all videos of a channel,
# so the per-video-group link rows are gone from the shared strings table
# and the channel name / link columns are rewired to the new consolidated
# channel links. The hyperlinks already attached to column C keep pointing
# at their original targets - only the displayed text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "sidemen shorts",
    "anthony padilla",
    "mr beast",
    "nile red",
    "ksi clips",
    "impalsive clips",
    "w2s clips",
    "smosh shorts"
)

$links = @(
    "https://www.youtube.com/c/UCbAZH3nTxzyNmehmTUhuUsA",
    "https://www.youtube.com/c/UCPJHQ5_DLtxZ1gzBvZE99_g",
    "https://www.youtube.com/channel/UC4-79UOlP48-QNGgCko5p2g",
    "https://www.youtube.com/c/UCA0mlN90EHCizvo101nbr-g",
    "https://www.youtube.com/channel/UCMiY4t431lhXY4QtPZtzftQ",
    "https://www.youtube.com/c/UCE9ZKI1b_PhVm3gejYuilhw",
    "https://www.youtube.com/c/UCZiJzk4wTIzaqHI4FXZ_eRQ",
    "https://www.youtube.com/c/UCS_NmOvbqaC9ccWSymx5Gpg"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $links[$i]
}

$ws.Range("C9").Select()
